$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 704.5
$ws.Range("I41").Value = 747.7143
$ws.Range("J41").Value = 402
$ws.Range("K41").Value = 747.7143
$ws.Range("L41").Value = 402
$ws.Range("M41").Value = -307.7143
$ws.Range("N41").Value = -1282

$ws.Range("H135").Value = 1217.7273
$ws.Range("I135").Value = 1036
$ws.Range("J135").Value = 3035
$ws.Range("K135").Value = 9324
$ws.Range("L135").Value = 27315
$ws.Range("M135").Value = -6789
$ws.Range("N135").Value = -32385

$ws.Range("H137").Value = 2245.4211
$ws.Range("I137").Value = 1341.2858
$ws.Range("J137").Value = 2772.8333
$ws.Range("K137").Value = 4023.8574
$ws.Range("L137").Value = 8318.499899999999
$ws.Range("M137").Value = -1473.8574
$ws.Range("N137").Value = -13418.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1095
$ws.Range("I2").Value = 1095
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1095
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -982
$ws.Range("N2").Value = $null

$ws.Range("H43").Value = 5028388
$ws.Range("I43").Value = 10010111
$ws.Range("J43").Value = 46664.668
$ws.Range("K43").Value = 10010111
$ws.Range("L43").Value = 46664.668
$ws.Range("M43").Value = -10009798
$ws.Range("N43").Value = -47290.668

$ws.Range("H61").Value = 2552.2222
$ws.Range("I61").Value = 1710.8572
$ws.Range("J61").Value = 5497
$ws.Range("K61").Value = 1710.8572
$ws.Range("L61").Value = 5497
$ws.Range("M61").Value = -1498.8572
$ws.Range("N61").Value = -5921

$ws.Range("H88").Value = 1409.6666
$ws.Range("I88").Value = 1457.3334
$ws.Range("J88").Value = 1266.6666
$ws.Range("K88").Value = 1457.3334
$ws.Range("L88").Value = 1266.6666
$ws.Range("M88").Value = -1051.3334
$ws.Range("N88").Value = -2078.6666

$ws.Range("H91").Value = 1409.6666
$ws.Range("I91").Value = 1457.3334
$ws.Range("J91").Value = 1266.6666
$ws.Range("K91").Value = 1457.3334
$ws.Range("L91").Value = 1266.6666
$ws.Range("M91").Value = -53.33339999999998
$ws.Range("N91").Value = -4074.6666

$ws.Range("H116").Value = 1095
$ws.Range("I116").Value = 1095
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1095
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1199
$ws.Range("N116").Value = $null

$ws.Range("H122").Value = 1696.8
$ws.Range("I122").Value = 1371.25
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 4113.75
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -1663.75
$ws.Range("N122").Value = -13897

$ws.Range("H124").Value = 27633
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 27633
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 27633
$ws.Range("N124").Value = -37453

$ws.Range("H136").Value = 2552.2222
$ws.Range("I136").Value = 1710.8572
$ws.Range("J136").Value = 5497
$ws.Range("K136").Value = 5132.571599999999
$ws.Range("L136").Value = 16491
$ws.Range("M136").Value = -2582.571599999999
$ws.Range("N136").Value = -21591

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1095
$ws.Range("I3").Value = 1095
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1095
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -981
$ws.Range("N3").Value = $null

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4907.9614
$ws.Range("I31").Value = 1826.125
$ws.Range("J31").Value = 6277.6665
$ws.Range("K31").Value = 1826.125
$ws.Range("L31").Value = 6277.6665
$ws.Range("M31").Value = -1531.125
$ws.Range("N31").Value = -6867.6665

$ws.Range("H34").Value = 4907.9614
$ws.Range("I34").Value = 1826.125
$ws.Range("J34").Value = 6277.6665
$ws.Range("K34").Value = 1826.125
$ws.Range("L34").Value = 6277.6665
$ws.Range("M34").Value = -1624.125
$ws.Range("N34").Value = -6681.6665

$ws.Range("H62").Value = 2520
$ws.Range("I62").Value = 2100
$ws.Range("J62").Value = 2940
$ws.Range("K62").Value = 2100
$ws.Range("L62").Value = 2940
$ws.Range("M62").Value = -1476
$ws.Range("N62").Value = -4188

$ws.Range("H65").Value = 2520
$ws.Range("I65").Value = 2100
$ws.Range("J65").Value = 2940
$ws.Range("K65").Value = 10500
$ws.Range("L65").Value = 14700
$ws.Range("M65").Value = -7380
$ws.Range("N65").Value = -20940

$ws.Range("H68").Value = 88882.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 88882.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 88882.5
$ws.Range("N68").Value = -90380.5

$ws.Range("H71").Value = 88882.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 88882.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 266647.5
$ws.Range("N71").Value = -274135.5

$ws.Range("H132").Value = 2105.9546
$ws.Range("I132").Value = 1920.5238
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 5761.5714
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -3231.5714
$ws.Range("N132").Value = -23060

$ws.Range("H134").Value = 2735.5334
$ws.Range("I134").Value = 2875.4285
$ws.Range("J134").Value = 777
$ws.Range("K134").Value = 8626.2855
$ws.Range("L134").Value = 2331
$ws.Range("M134").Value = -6091.2855
$ws.Range("N134").Value = -7401

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2063.2666
$ws.Range("I131").Value = 1224.8334
$ws.Range("J131").Value = 2622.2222
$ws.Range("K131").Value = 3674.5002
$ws.Range("L131").Value = 7866.6666
$ws.Range("M131").Value = 1365.4998
$ws.Range("N131").Value = -17946.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 855.5
$ws.Range("I22").Value = 468.8

$ws.Range("H27").Value = 855.5
$ws.Range("I27").Value = 468.8

$ws.Range("H44").Value = 4000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 4000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 4000
$ws.Range("N44").Value = -4912

$ws.Range("H63").Value = 41249.668
$ws.Range("I63").Value = 41249.668
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 41249.668
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -40500.668

$ws.Range("H66").Value = 41249.668
$ws.Range("I66").Value = 41249.668
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 123749.004
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -120005.004

$ws.Range("H68").Value = 7560.3335
$ws.Range("I68").Value = 6135.75
$ws.Range("J68").Value = 8700
$ws.Range("K68").Value = 6135.75
$ws.Range("L68").Value = 8700
$ws.Range("M68").Value = -5386.75
$ws.Range("N68").Value = -10198

$ws.Range("H71").Value = 7560.3335
$ws.Range("I71").Value = 6135.75
$ws.Range("J71").Value = 8700
$ws.Range("K71").Value = 30678.75
$ws.Range("L71").Value = 43500
$ws.Range("M71").Value = -26934.75
$ws.Range("N71").Value = -50988

$ws.Range("H132").Value = 3448.6667
$ws.Range("I132").Value = 3674.25
$ws.Range("J132").Value = 2997.5
$ws.Range("K132").Value = 11022.75
$ws.Range("L132").Value = 8992.5
$ws.Range("M132").Value = -8492.75
$ws.Range("N132").Value = -14052.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 14301
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 14301
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 14301
$ws.Range("N82").Value = -15067

$ws.Range("H85").Value = 14301
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 14301
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 14301
$ws.Range("N85").Value = -16953

$ws.Range("H96").Value = 1250
$ws.Range("I96").Value = 1250
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1250
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 123
$ws.Range("N96").Value = $null

$ws.Range("H126").Value = 4419.857
$ws.Range("I126").Value = 1699
$ws.Range("J126").Value = 7140.7144
$ws.Range("K126").Value = 5097
$ws.Range("L126").Value = 21422.1432
$ws.Range("M126").Value = -2627
$ws.Range("N126").Value = -26362.1432

$ws.Range("H132").Value = 1259.8182
$ws.Range("I132").Value = 1151.8572
$ws.Range("J132").Value = 1448.75
$ws.Range("K132").Value = 3455.5716
$ws.Range("L132").Value = 4346.25
$ws.Range("M132").Value = -925.5715999999998
$ws.Range("N132").Value = -9406.25

$ws.Range("H136").Value = 3548.2856
$ws.Range("I136").Value = 2479.75
$ws.Range("J136").Value = 4973
$ws.Range("K136").Value = 7439.25
$ws.Range("L136").Value = 14919
$ws.Range("M136").Value = -4889.25
$ws.Range("N136").Value = -20019
